# update new file(Database) to solve the conflict
#
# Sets a cell's value as literal text (shared-string), even when the text
# looks like a number/date, while leaving the cell's style untouched.
# Plain `Range.Value = "123"` gets auto-coerced to a number by Excel, and
# forcing text via NumberFormat="@" leaves a residual "quote prefix" style
# behind -- so we clear formats after the write and, if the cell needs to
# keep its original style, repaint that style from an untouched donor cell
# that already carries it (PasteSpecial formats only; doesn't touch value).
function Set-TextCell {
    param($cell, [string]$text, $donor = $null)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
    if ($donor -ne $null) {
        $donor.Copy()
        $cell.PasteSpecial(-4122)
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Customer sheet: append two new customers (rows 6 and 7)
# ---------------------------------------------------------------------
$customer = $wb.Worksheets.Item("Customer")

Set-TextCell $customer.Cells.Item(6,1) "5"
Set-TextCell $customer.Cells.Item(6,2) "yolanda"
Set-TextCell $customer.Cells.Item(6,3) "FEMALE"
Set-TextCell $customer.Cells.Item(6,4) "1"
Set-TextCell $customer.Cells.Item(6,5) "10"
Set-TextCell $customer.Cells.Item(6,6) "123456"
Set-TextCell $customer.Cells.Item(6,7) "12332112345"

Set-TextCell $customer.Cells.Item(7,1) "6"
Set-TextCell $customer.Cells.Item(7,2) "ddd"
Set-TextCell $customer.Cells.Item(7,3) "FEMALE"
Set-TextCell $customer.Cells.Item(7,4) "0"
Set-TextCell $customer.Cells.Item(7,5) "1"
Set-TextCell $customer.Cells.Item(7,6) "234567"
Set-TextCell $customer.Cells.Item(7,7) "12345678901"

# ---------------------------------------------------------------------
# Drink sheet: drop the test "yolanda" drink row, fix up stock/sell counts
# ---------------------------------------------------------------------
$drink = $wb.Worksheets.Item("Drink")

$drink.Rows.Item(6).Delete()

Set-TextCell $drink.Cells.Item(5,4) "999"
Set-TextCell $drink.Cells.Item(5,5) "1"

# ---------------------------------------------------------------------
# Food sheet: fix up stock/sell counts for rows 2 and 3
# ---------------------------------------------------------------------
$food = $wb.Worksheets.Item("Food")

$foodD1 = $food.Cells.Item(1,4)
$foodE1 = $food.Cells.Item(1,5)
$foodA2 = $food.Cells.Item(2,1)
$foodC2 = $food.Cells.Item(2,3)

Set-TextCell $food.Cells.Item(2,4) "999" $foodD1
Set-TextCell $food.Cells.Item(2,5) "13" $foodE1

Set-TextCell $food.Cells.Item(3,1) "2" $foodA2
Set-TextCell $food.Cells.Item(3,3) "20.0" $foodC2
Set-TextCell $food.Cells.Item(3,4) "139" $foodD1
Set-TextCell $food.Cells.Item(3,5) "61" $foodE1

# ---------------------------------------------------------------------
# Room sheet: update the first booking row's values
# ---------------------------------------------------------------------
$room = $wb.Worksheets.Item("Room")

$roomA3 = $room.Cells.Item(3,1)
$roomC3 = $room.Cells.Item(3,3)
$roomD3 = $room.Cells.Item(3,4)
$roomE3 = $room.Cells.Item(3,5)
$roomF3 = $room.Cells.Item(3,6)

Set-TextCell $room.Cells.Item(2,1) "1" $roomA3
Set-TextCell $room.Cells.Item(2,3) "188.0" $roomC3
Set-TextCell $room.Cells.Item(2,4) "2022-12-11" $roomD3
Set-TextCell $room.Cells.Item(2,5) "5" $roomE3
Set-TextCell $room.Cells.Item(2,6) "10" $roomF3
